$wb = $excel.ActiveWorkbook

$wsEco = $wb.Worksheets.Item("ecological_params")
$wsScale = $wb.Worksheets.Item("Scaling")

# Rename the PP_cuboid_* headers to their shorter names.
$wsEco.Range("B1").Value = "PP_virgin"
$wsEco.Range("C1").Value = "PP_recycled"
$wsEco.Range("D1").Value = "PP_recycled_industrial"

# The ecological_params tab becomes the active/selected tab again.
$wsEco.Activate()
